# Updates the ozone SHAP feature-importance table (Sheet1, columns A:C, rows 2-106)
# with refreshed SHAP values. Because the sheet is kept sorted in descending
# order of the "shap" column, re-running the SHAP computation also changes the
# row order for the "feat" index (A) / label (B) columns further down the
# table even though the underlying feature set (0-104) is unchanged.
#
# We therefore rewrite the full A2:C106 block in one shot via a bulk array
# assignment rather than patching individual cells, which both matches the
# target state exactly and is far fewer COM round-trips than 300+ single
# cell writes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 105,3
$arr[0,0] = 74
$arr[0,1] = '$\langle qq \vert qq \rangle$'
$arr[0,2] = [double]"0.00194509764621996"
$arr[1,0] = 12
$arr[1,1] = 'h$_{q}$'
$arr[1,2] = [double]"0.001806405169951488"
$arr[2,0] = 76
$arr[2,1] = '$\langle ss \vert ss \rangle$'
$arr[2,2] = [double]"0.001030939165310405"
$arr[3,0] = 37
$arr[3,1] = '$F_{q}$'
$arr[3,2] = [double]"0.0008744899204508957"
$arr[4,0] = 101
$arr[4,1] = '$(\langle pq \vert pq \rangle)_{3}$'
$arr[4,2] = [double]"0.0006352018862609988"
$arr[5,0] = 35
$arr[5,1] = '$F_{q}^{\text{SCF}}$'
$arr[5,2] = [double]"0.0006286188149934027"
$arr[6,0] = 22
$arr[6,1] = 'h$_{s}$'
$arr[6,2] = [double]"0.0004868368230681382"
$arr[7,0] = 99
$arr[7,1] = '$(\langle pp \vert pp \rangle)_{3}$'
$arr[7,2] = [double]"0.0003313962474878542"
$arr[8,0] = 63
$arr[8,1] = '$(F_{p}^{\text{SCF}})_{3}$'
$arr[8,2] = [double]"0.0003044498699009141"
$arr[9,0] = 3
$arr[9,1] = 'h$_{p}^{3}$'
$arr[9,2] = [double]"0.0002683088222694"
$arr[10,0] = 43
$arr[10,1] = '$F_{s}^{\text{SCF}}$'
$arr[10,2] = [double]"0.0002593594698316075"
$arr[11,0] = 45
$arr[11,1] = '$F_{s}$'
$arr[11,2] = [double]"0.000252806681707469"
$arr[12,0] = 21
$arr[12,1] = 'h$_{rs}^{3}$'
$arr[12,2] = [double]"0.0002470983937014804"
$arr[13,0] = 102
$arr[13,1] = '$(\langle pq \vert qp \rangle)_{3}$'
$arr[13,2] = [double]"0.000155506172445445"
$arr[14,0] = 30
$arr[14,1] = 'FA$_{qs}$'
$arr[14,2] = [double]"0.0001549540278802658"
$arr[15,0] = 57
$arr[15,1] = '$(F_{p})_{2}$'
$arr[15,2] = [double]"0.0001427543223523081"
$arr[16,0] = 29
$arr[16,1] = 'FI$_{qs}$'
$arr[16,2] = [double]"0.0001406777453182349"
$arr[17,0] = 67
$arr[17,1] = '$(F_{r}^{\text{SCF}})_{3}$'
$arr[17,2] = [double]"0.0001132686806636267"
$arr[18,0] = 65
$arr[18,1] = '$(F_{p})_{3}$'
$arr[18,2] = [double]"0.000103728263892844"
$arr[19,0] = 104
$arr[19,1] = '$(\langle rs \vert sr \rangle)_{3}$'
$arr[19,2] = [double]"0.000100844566747294"
$arr[20,0] = 2
$arr[20,1] = 'h$_{p}^{2}$'
$arr[20,2] = [double]"9.902569680403246e-05"
$arr[21,0] = 91
$arr[21,1] = '$(\langle pp \vert pp \rangle)_{2}$'
$arr[21,2] = [double]"9.683742161502815e-05"
$arr[22,0] = 103
$arr[22,1] = '$(\langle rs\vert rs \rangle)_{3}$'
$arr[22,2] = [double]"8.935739885445854e-05"
$arr[23,0] = 17
$arr[23,1] = 'h$_{r}^{3}$'
$arr[23,2] = [double]"8.793148476700385e-05"
$arr[24,0] = 93
$arr[24,1] = '$(\langle pq \vert pq \rangle)_{2}$'
$arr[24,2] = [double]"8.731045049528486e-05"
$arr[25,0] = 42
$arr[25,1] = '$(\eta_{r})_{0}$'
$arr[25,2] = [double]"8.729929690544922e-05"
$arr[26,0] = 15
$arr[26,1] = 'h$_{r}^{1}$'
$arr[26,2] = [double]"8.508264073544888e-05"
$arr[27,0] = 7
$arr[27,1] = 'h$_{pq}^{3}$'
$arr[27,2] = [double]"8.482727735757262e-05"
$arr[28,0] = 13
$arr[28,1] = 'h$_{qs}$'
$arr[28,2] = [double]"8.304035830667969e-05"
$arr[29,0] = 0
$arr[29,1] = 'h$_{p}^{0}$'
$arr[29,2] = [double]"8.144388374488473e-05"
$arr[30,0] = 26
$arr[30,1] = 'typ_3'
$arr[30,2] = [double]"5.77602566001905e-05"
$arr[31,0] = 100
$arr[31,1] = '$(\langle rr \vert rr \rangle)_{3}$'
$arr[31,2] = [double]"5.601211829127773e-05"
$arr[32,0] = 47
$arr[32,1] = '$(F_{p}^{\text{SCF}})_{1}$'
$arr[32,2] = [double]"5.01981901654735e-05"
$arr[33,0] = 83
$arr[33,1] = '$(\langle pp \vert pp \rangle)_{1}$'
$arr[33,2] = [double]"4.638268573879573e-05"
$arr[34,0] = 46
$arr[34,1] = '$\eta_{s}$'
$arr[34,2] = [double]"4.311554710674413e-05"
$arr[35,0] = 97
$arr[35,1] = '$(\langle pq \vert rs \rangle)_{3}$'
$arr[35,2] = [double]"4.070445902978927e-05"
$arr[36,0] = 9
$arr[36,1] = 'h$_{pr}^{1}$'
$arr[36,2] = [double]"3.885330499007429e-05"
$arr[37,0] = 85
$arr[37,1] = '$(\langle pq \vert pq \rangle)_{1}$'
$arr[37,2] = [double]"3.784229247006348e-05"
$arr[38,0] = 77
$arr[38,1] = '$(\langle pq \vert pq \rangle)_{0}$'
$arr[38,2] = [double]"3.720844953883601e-05"
$arr[39,0] = 20
$arr[39,1] = 'h$_{rs}^{2}$'
$arr[39,2] = [double]"3.717815880197529e-05"
$arr[40,0] = 16
$arr[40,1] = 'h$_{r}^{2}$'
$arr[40,2] = [double]"3.700131861588328e-05"
$arr[41,0] = 10
$arr[41,1] = 'h$_{pr}^{2}$'
$arr[41,2] = [double]"3.504774087439861e-05"
$arr[42,0] = 95
$arr[42,1] = '$(\langle rs\vert rs \rangle)_{2}$'
$arr[42,2] = [double]"3.370782852659465e-05"
$arr[43,0] = 5
$arr[43,1] = 'h$_{pq}^{1}$'
$arr[43,2] = [double]"3.232924699069303e-05"
$arr[44,0] = 23
$arr[44,1] = 'typ_0'
$arr[44,2] = [double]"3.209007779275085e-05"
$arr[45,0] = 94
$arr[45,1] = '$(\langle pq \vert qp \rangle)_{2}$'
$arr[45,2] = [double]"3.203533073389536e-05"
$arr[46,0] = 96
$arr[46,1] = '$(\langle rs \vert sr \rangle)_{2}$'
$arr[46,2] = [double]"3.164213345520241e-05"
$arr[47,0] = 55
$arr[47,1] = '$(F_{p}^{\text{SCF}})_{2}$'
$arr[47,2] = [double]"3.034945791054597e-05"
$arr[48,0] = 92
$arr[48,1] = '$(\langle rr \vert rr \rangle)_{2}$'
$arr[48,2] = [double]"2.943404609335781e-05"
$arr[49,0] = 6
$arr[49,1] = 'h$_{pq}^{2}$'
$arr[49,2] = [double]"2.678073466255247e-05"
$arr[50,0] = 89
$arr[50,1] = '$(\langle pq \vert rs \rangle)_{2}$'
$arr[50,2] = [double]"2.366151817830605e-05"
$arr[51,0] = 80
$arr[51,1] = '$(\langle rs \vert sr \rangle)_{0}$'
$arr[51,2] = [double]"2.202977823572696e-05"
$arr[52,0] = 87
$arr[52,1] = '$(\langle rs\vert rs \rangle)_{1}$'
$arr[52,2] = [double]"2.17615930497546e-05"
$arr[53,0] = 78
$arr[53,1] = '$(\langle pq \vert qp \rangle)_{0}$'
$arr[53,2] = [double]"2.144556018481902e-05"
$arr[54,0] = 86
$arr[54,1] = '$(\langle pq \vert qp \rangle)_{1}$'
$arr[54,2] = [double]"2.053014813978082e-05"
$arr[55,0] = 11
$arr[55,1] = 'h$_{pr}^{3}$'
$arr[55,2] = [double]"2.047007399487397e-05"
$arr[56,0] = 59
$arr[56,1] = '$(F_{r}^{\text{SCF}})_{2}$'
$arr[56,2] = [double]"1.931317148690808e-05"
$arr[57,0] = 51
$arr[57,1] = '$(F_{r}^{\text{SCF}})_{1}$'
$arr[57,2] = [double]"1.862987545932772e-05"
$arr[58,0] = 75
$arr[58,1] = '$(\langle rr \vert rr \rangle)_{0}$'
$arr[58,2] = [double]"1.777275067913309e-05"
$arr[59,0] = 71
$arr[59,1] = '$(\langle pq \vert rs \rangle)_{0}$'
$arr[59,2] = [double]"1.704935277646223e-05"
$arr[60,0] = 14
$arr[60,1] = 'h$_{r}^{0}$'
$arr[60,2] = [double]"1.609666314417677e-05"
$arr[61,0] = 69
$arr[61,1] = '$(F_{r})_{3}$'
$arr[61,2] = [double]"1.49898674812587e-05"
$arr[62,0] = 18
$arr[62,1] = 'h$_{rs}^{0}$'
$arr[62,2] = [double]"1.495013778519909e-05"
$arr[63,0] = 4
$arr[63,1] = 'h$_{pq}^{0}$'
$arr[63,2] = [double]"1.453603244127391e-05"
$arr[64,0] = 1
$arr[64,1] = 'h$_{p}^{1}$'
$arr[64,2] = [double]"1.416397742143847e-05"
$arr[65,0] = 25
$arr[65,1] = 'typ_2'
$arr[65,2] = [double]"1.348391260600136e-05"
$arr[66,0] = 73
$arr[66,1] = '$(\langle pp \vert pp \rangle)_{0}$'
$arr[66,2] = [double]"1.308480641036284e-05"
$arr[67,0] = 81
$arr[67,1] = '$(\langle pq \vert rs \rangle)_{1}$'
$arr[67,2] = [double]"1.156695138761045e-05"
$arr[68,0] = 8
$arr[68,1] = 'h$_{pr}^{0}$'
$arr[68,2] = [double]"1.12419102293991e-05"
$arr[69,0] = 31
$arr[69,1] = '$(F_{p}^{\text{SCF}})_{0}$'
$arr[69,2] = [double]"1.119905882560219e-05"
$arr[70,0] = 79
$arr[70,1] = '$(\langle rs\vert rs \rangle)_{0}$'
$arr[70,2] = [double]"1.1162357291741e-05"
$arr[71,0] = 24
$arr[71,1] = 'typ_1'
$arr[71,2] = [double]"9.381280056516233e-06"
$arr[72,0] = 61
$arr[72,1] = '$(F_{r})_{2}$'
$arr[72,2] = [double]"9.124969217982509e-06"
$arr[73,0] = 53
$arr[73,1] = '$(F_{r})_{1}$'
$arr[73,2] = [double]"9.038984198279314e-06"
$arr[74,0] = 41
$arr[74,1] = '$(F_{r})_{0}$'
$arr[74,2] = [double]"8.733291539230997e-06"
$arr[75,0] = 62
$arr[75,1] = '$(\eta_{r})_{2}$'
$arr[75,2] = [double]"7.739212165113958e-06"
$arr[76,0] = 19
$arr[76,1] = 'h$_{rs}^{1}$'
$arr[76,2] = [double]"7.63096311561106e-06"
$arr[77,0] = 88
$arr[77,1] = '$(\langle rs \vert sr \rangle)_{1}$'
$arr[77,2] = [double]"7.139091976137434e-06"
$arr[78,0] = 84
$arr[78,1] = '$(\langle rr \vert rr \rangle)_{1}$'
$arr[78,2] = [double]"7.08190507842165e-06"
$arr[79,0] = 50
$arr[79,1] = '$(\eta_{p})_{1}$'
$arr[79,2] = [double]"5.761066374722253e-06"
$arr[80,0] = 34
$arr[80,1] = '$(\eta_{p})_{0}$'
$arr[80,2] = [double]"4.56872119437079e-06"
$arr[81,0] = 49
$arr[81,1] = '$(F_{p})_{1}$'
$arr[81,2] = [double]"4.472066202883751e-06"
$arr[82,0] = 39
$arr[82,1] = '$(F_{r}^{\text{SCF}})_{0}$'
$arr[82,2] = [double]"3.738925324421023e-06"
$arr[83,0] = 58
$arr[83,1] = '$(\eta_{p})_{2}$'
$arr[83,2] = [double]"2.902489172849546e-06"
$arr[84,0] = 33
$arr[84,1] = '$(F_{p})_{0}$'
$arr[84,2] = [double]"5.551928990137636e-07"
$arr[85,0] = 38
$arr[85,1] = '$\eta_{q}$'
$arr[85,2] = [double]"4.477981570927126e-07"
$arr[86,0] = 70
$arr[86,1] = '$(\eta_{r})_{3}$'
$arr[86,2] = [double]"4.375728308215206e-07"
$arr[87,0] = 66
$arr[87,1] = '$(\eta_{p})_{3}$'
$arr[87,2] = [double]"3.692255639771814e-07"
$arr[88,0] = 98
$arr[88,1] = '$(\langle pq \vert sr \rangle)_{3}$'
$arr[88,2] = [double]"3.505646844630095e-07"
$arr[89,0] = 90
$arr[89,1] = '$(\langle pq \vert sr \rangle)_{2}$'
$arr[89,2] = [double]"2.682307745673009e-07"
$arr[90,0] = 72
$arr[90,1] = '$(\langle pq \vert sr \rangle)_{0}$'
$arr[90,2] = [double]"2.502153056358304e-07"
$arr[91,0] = 82
$arr[91,1] = '$(\langle pq \vert sr \rangle)_{1}$'
$arr[91,2] = [double]"2.121715547427413e-07"
$arr[92,0] = 40
$arr[92,1] = '$(\omega_{r})_{0}$'
$arr[92,2] = [double]"1.6989203585197e-07"
$arr[93,0] = 48
$arr[93,1] = '$(\omega_{p})_{1}$'
$arr[93,2] = [double]"1.591946705930024e-07"
$arr[94,0] = 32
$arr[94,1] = '$(\omega_{p})_{0}$'
$arr[94,2] = [double]"1.379493345633971e-07"
$arr[95,0] = 56
$arr[95,1] = '$(\omega_{p})_{2}$'
$arr[95,2] = [double]"1.316511079805008e-07"
$arr[96,0] = 54
$arr[96,1] = '$(\eta_{r})_{1}$'
$arr[96,2] = [double]"1.24971359396434e-07"
$arr[97,0] = 52
$arr[97,1] = '$(\omega_{r})_{1}$'
$arr[97,2] = [double]"1.121639238941802e-07"
$arr[98,0] = 44
$arr[98,1] = '$\omega_{s}$'
$arr[98,2] = [double]"1.002623803145276e-07"
$arr[99,0] = 64
$arr[99,1] = '$(\omega_{p})_{3}$'
$arr[99,2] = [double]"8.850376789605033e-08"
$arr[100,0] = 68
$arr[100,1] = '$(\omega_{r})_{3}$'
$arr[100,2] = [double]"8.112541656118031e-08"
$arr[101,0] = 60
$arr[101,1] = '$(\omega_{r})_{2}$'
$arr[101,2] = [double]"6.269400613764571e-08"
$arr[102,0] = 28
$arr[102,1] = 'F$_{qs}$'
$arr[102,2] = [double]"1.292334876238133e-08"
$arr[103,0] = 27
$arr[103,1] = '$\mathbf{b}$'
$arr[103,2] = [double]"1.169643555158222e-08"
$arr[104,0] = 36
$arr[104,1] = '$\omega_{q}$'
$arr[104,2] = [double]"1.007245330045063e-08"

$ws.Range("A2:C106").Value = $arr

Write-Host "Updated $($arr.GetLength(0)) SHAP rows (A2:C106)"
